$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the Price/Volume columns as text first so that numeric-looking
# strings (e.g. "7.10", "4.40") are preserved exactly instead of being
# reinterpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '68.377.37'
$ws.Range("E2").Value = '  -1.87%  '
$ws.Range("D3").Value = '2.451.36'
$ws.Range("E3").Value = '  -1.87%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '557.25'
$ws.Range("E5").Value = '  -2.24%  '
$ws.Range("D6").Value = '162.95'
$ws.Range("E6").Value = '  -1.77%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '0.501'
$ws.Range("E8").Value = '  -2.13%  '
$ws.Range("D9").Value = '2.450.02'
$ws.Range("E9").Value = '  -1.88%  '
$ws.Range("E10").Value = '  -5.95%  '
$ws.Range("E11").Value = '  -1.86%  '
$ws.Range("E12").Value = '  -6.02%  '
$ws.Range("D13").Value = '4.76'
$ws.Range("E13").Value = '  -3.32%  '
$ws.Range("D14").Value = '2.901.74'
$ws.Range("E14").Value = '  -1.67%  '
$ws.Range("D15").Value = '68.322.90'
$ws.Range("D16").Value = '0.0000167'
$ws.Range("E16").Value = '  -4.46%  '
$ws.Range("D17").Value = '23.11'
$ws.Range("E17").Value = '  -5.23%  '
$ws.Range("D18").Value = '2.463.56'
$ws.Range("E18").Value = '  -1.04%  '
$ws.Range("D19").Value = '10.83'
$ws.Range("E19").Value = '  -3.30%  '
$ws.Range("D20").Value = '7.11'
$ws.Range("E20").Value = '  -3.85%  '
$ws.Range("D21").Value = '339.52'
$ws.Range("E21").Value = '  -2.13%  '
$ws.Range("D22").Value = '3.74'
$ws.Range("E22").Value = '  -3.48%  '
$ws.Range("E23").Value = '  -0.17%  '
$ws.Range("E24").Value = '  -4.67%  '
$ws.Range("D25").Value = '67.22'
$ws.Range("E25").Value = '  -4.59%  '
$ws.Range("D26").Value = '3.65'
$ws.Range("E26").Value = '  -6.55%  '
$ws.Range("D27").Value = '2.577.74'
$ws.Range("E27").Value = '  -1.44%  '
$ws.Range("D28").Value = '1.01'
$ws.Range("E28").Value = '  +0.85%  '
$ws.Range("D29").Value = '8.04'
$ws.Range("E29").Value = '  -7.33%  '
$ws.Range("E30").Value = '  -6.90%  '
$ws.Range("D31").Value = '7.10'
$ws.Range("E31").Value = '  -9.34%  '
$ws.Range("E32").Value = '  +0.05%  '
$ws.Range("D33").Value = '423.44'
$ws.Range("E33").Value = '  -4.95%  '
$ws.Range("E34").Value = '  -4.23%  '
$ws.Range("E35").Value = '  -4.41%  '
$ws.Range("D36").Value = '157.39'
$ws.Range("E36").Value = '  +1.06%  '
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("E39").Value = '  -4.86%  '
$ws.Range("D40").Value = '17.76'
$ws.Range("E40").Value = '  -2.48%  '
$ws.Range("D41").Value = '0.302'
$ws.Range("E41").Value = '  -4.25%  '
$ws.Range("D42").Value = '4.40'
$ws.Range("E42").Value = '  -4.80%  '
$ws.Range("D43").Value = '1.49'
$ws.Range("E43").Value = '  -6.04%  '
$ws.Range("E44").Value = '  +0.26%  '
$ws.Range("D45").Value = '133.68'
$ws.Range("E45").Value = '  -4.54%  '
$ws.Range("E46").Value = '  -6.78%  '
$ws.Range("E47").Value = '  -4.04%  '
$ws.Range("D48").Value = '0.0714'
$ws.Range("E48").Value = '  -2.14%  '
$ws.Range("D49").Value = '0.476'
$ws.Range("E49").Value = '  -7.44%  '
$ws.Range("D50").Value = '0.560'
$ws.Range("E50").Value = '  -2.62%  '
$ws.Range("D51").Value = '0.0904'
$ws.Range("E51").Value = '  -2.13%  '

# Restore the default cell style so no stray number-format styling is left behind.
$ws.Range("D2:E51").Style = "Normal"

